# [Pipette] Battery Connector 추가 - solder cap type
$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> "Battery CON"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Battery CON"

# Column widths (Excel best-fit for column B content)
$ws2.Columns.Item(2).ColumnWidth = 11

# Row 5 (string cells first so shared-string indices line up with the target file: 9155-000, AVX)
$ws2.Range("B5").Value = "9155-000"
$ws2.Range("C5").Value = "AVX"

# Header row (row 4): Body, Spring, Total
$ws2.Range("D4").Value = "Body"
$ws2.Range("E4").Value = "Spring"
$ws2.Range("F4").Value = "Total"

# Row 6: K113868632, Top-Link
$ws2.Range("B6").Value = "K113868632"
$ws2.Range("C6").Value = "Top-Link"

# Row 7: Tyco
$ws2.Range("C7").Value = "Tyco"

# Numeric / formula cells
$ws2.Range("D5").Value = 3
$ws2.Range("E5").Value = 2.3
$ws2.Range("F5").Formula = "=D5+E5"
$ws2.Range("F6").Value = 5.3
$ws2.Range("F7").Value = 6.2

# Row 8
$ws2.Range("D8").Value = 3.45
$ws2.Range("E8").Formula = "=4.9-3.45"
$ws2.Range("F8").Formula = "=D8+E8"

# Row 12-15
$ws2.Range("O12").Value = 0.089
$ws2.Range("O13").Value = 0.204
$ws2.Range("Q13").Value = 0.182
$ws2.Range("O14").Formula = "=O13-O12"
$ws2.Range("Q14").Value = 0.11
$ws2.Range("Q15").Formula = "=Q13-Q14"

# Selection on Battery CON sheet
$ws2.Range("L24").Select()

# Activate Battery CON tab (makes it tabSelected + workbook activeTab)
$ws2.Activate()

$wb.Save()
